$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H17").Value = "https://www.auchan.fr/lego-star-wars75417-le-marcheur-at-st/pr-C1844702"
$ws.Range("J17").Value = "https://www.carrefour.fr/p/lego-le-marcheur-at-st-75417-lego-5702017817668"
